$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Real SP Sprint 4" (column J) values for rows 20-26 (skip row 23, already filled)
$ws.Range("J20").Value = 46
$ws.Range("J21").Value = 70
$ws.Range("J22").Value = 30
$ws.Range("J24").Value = 36
$ws.Range("J25").Value = 58
$ws.Range("J26").Value = 50

# Add the missing total formula for column J on the totals row
$ws.Range("J27").Formula = "=SUM(J20:J26)"

# Update the active selection to match the author's final cursor position
$ws.Range("F15").Select()
